$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 542
$ws.Range("I12").Value = 313.5
$ws.Range("K12").Value = 313.5
$ws.Range("M12").Value = -143.5

$ws.Range("H15").Value = 544916
$ws.Range("I15").Value = 544916
$ws.Range("K15").Value = 1634748
$ws.Range("M15").Value = -1634579

$ws.Range("H17").Value = 2524.2144
$ws.Range("J17").Value = 2524.2144
$ws.Range("L17").Value = 7572.6432
$ws.Range("N17").Value = -7908.6432

$ws.Range("H19").Value = 991.1539
$ws.Range("I19").Value = 951.7368
$ws.Range("J19").Value = 1098.1428
$ws.Range("K19").Value = 951.7368
$ws.Range("L19").Value = 1098.1428
$ws.Range("M19").Value = -776.7368
$ws.Range("N19").Value = -1448.1428

$ws.Range("H38").Value = 4889.9473
$ws.Range("I38").Value = 1702.8889
$ws.Range("K38").Value = 5108.6667
$ws.Range("M38").Value = -4736.6667

$ws.Range("H55").Value = 208.7
$ws.Range("I55").Value = 152.25
$ws.Range("J55").Value = 434.5
$ws.Range("K55").Value = 152.25
$ws.Range("L55").Value = 434.5
$ws.Range("M55").Value = 61.75
$ws.Range("N55").Value = -862.5

$ws.Range("H69").Value = 18087.176
$ws.Range("J69").Value = 21038.076
$ws.Range("L69").Value = 63114.228
$ws.Range("N69").Value = -64862.228

$ws.Range("H72").Value = 18087.176
$ws.Range("J72").Value = 21038.076
$ws.Range("L72").Value = 189342.684
$ws.Range("N72").Value = -198078.684

$ws.Range("H80").Value = 283.57144
$ws.Range("I80").Value = 190.83333
$ws.Range("J80").Value = 353.125
$ws.Range("K80").Value = 572.49999
$ws.Range("L80").Value = 1059.375
$ws.Range("M80").Value = 425.50001
$ws.Range("N80").Value = -3055.375

$ws.Range("H83").Value = 283.57144
$ws.Range("I83").Value = 190.83333
$ws.Range("J83").Value = 353.125
$ws.Range("K83").Value = 1717.49997
$ws.Range("L83").Value = 3178.125
$ws.Range("M83").Value = 3274.50003
$ws.Range("N83").Value = -13162.125

$ws.Range("H96").Value = 127.5
$ws.Range("J96").Value = 129.5
$ws.Range("L96").Value = 388.5
$ws.Range("N96").Value = -3134.5

$ws.Range("H132").Value = 2567.0938
$ws.Range("J132").Value = 2995
$ws.Range("L132").Value = 8985
$ws.Range("N132").Value = -14045

$ws.Range("H138").Value = 5258.4907
$ws.Range("I138").Value = 4576.516
$ws.Range("J138").Value = 6139.375
$ws.Range("K138").Value = 13729.548
$ws.Range("L138").Value = 18418.125
$ws.Range("M138").Value = -8589.547999999999
$ws.Range("N138").Value = -28698.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 62000
$ws.Range("I123").Value = 50000
$ws.Range("K123").Value = 50000
$ws.Range("M123").Value = -45100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3738.389
$ws.Range("I86").Value = 4068.7693
$ws.Range("K86").Value = 4068.7693
$ws.Range("M86").Value = -2945.7693

$ws.Range("H89").Value = 3738.389
$ws.Range("I89").Value = 4068.7693
$ws.Range("K89").Value = 20343.8465
$ws.Range("M89").Value = -14727.8465

$ws.Range("H94").Value = 2353.3125
$ws.Range("I94").Value = 2233.6667
$ws.Range("K94").Value = 2233.6667
$ws.Range("M94").Value = -1782.6667

$ws.Range("H105").Value = 4151.625
$ws.Range("I105").Value = 4069.6155
$ws.Range("K105").Value = 4069.6155
$ws.Range("M105").Value = -2322.6155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 20998.666
$ws.Range("J55").Value = 20998.666
$ws.Range("L55").Value = 20998.666
$ws.Range("N55").Value = -21628.666

$ws.Range("H62").Value = 4899.533
$ws.Range("I62").Value = 4443.8887
$ws.Range("J62").Value = 5583
$ws.Range("K62").Value = 4443.8887
$ws.Range("L62").Value = 5583
$ws.Range("M62").Value = -3819.8887
$ws.Range("N62").Value = -6831

$ws.Range("H65").Value = 4899.533
$ws.Range("I65").Value = 4443.8887
$ws.Range("J65").Value = 5583
$ws.Range("K65").Value = 22219.4435
$ws.Range("L65").Value = 27915
$ws.Range("M65").Value = -19099.4435
$ws.Range("N65").Value = -34155

$ws.Range("H94").Value = 5936.2144
$ws.Range("J94").Value = 5888.857
$ws.Range("L94").Value = 5888.857
$ws.Range("N94").Value = -6790.857

$ws.Range("H107").Value = 2220.25
$ws.Range("I107").Value = 2108.125
$ws.Range("J107").Value = 2444.5
$ws.Range("K107").Value = 2108.125
$ws.Range("L107").Value = 2444.5
$ws.Range("M107").Value = -188.125
$ws.Range("N107").Value = -6284.5

$ws.Range("H124").Value = 7044642
$ws.Range("J124").Value = 7044642
$ws.Range("L124").Value = 7044642
$ws.Range("N124").Value = -7049552

$ws.Range("H131").Value = 99999
$ws.Range("J131").Value = 99999
$ws.Range("L131").Value = 99999
$ws.Range("M131").Value = -110079

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 203.77777
$ws.Range("I2").Value = 146.38461
$ws.Range("K2").Value = 878.3076600000001
$ws.Range("M2").Value = -765.3076600000001

$ws.Range("H7").Value = 60.5
$ws.Range("I7").Value = 60.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 181.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -69.5
$ws.Range("N7").ClearContents()

$ws.Range("H64").Value = 6295.375
$ws.Range("I64").Value = 7999.857
$ws.Range("J64").Value = 4969.6665
$ws.Range("K64").Value = 23999.571
$ws.Range("L64").Value = 14908.9995
$ws.Range("M64").Value = -23729.571
$ws.Range("N64").Value = -15448.9995

$ws.Range("H67").Value = 6295.375
$ws.Range("I67").Value = 7999.857
$ws.Range("J67").Value = 4969.6665
$ws.Range("K67").Value = 23999.571
$ws.Range("L67").Value = 14908.9995
$ws.Range("M67").Value = -23063.571
$ws.Range("N67").Value = -16780.9995

$ws.Range("H68").Value = 846.4
$ws.Range("J68").Value = 1400.6666
$ws.Range("L68").Value = 4201.9998
$ws.Range("N68").Value = -5823.9998

$ws.Range("H71").Value = 846.4
$ws.Range("J71").Value = 1400.6666
$ws.Range("L71").Value = 12605.9994
$ws.Range("N71").Value = -20717.9994

$ws.Range("H122").Value = 348.5
$ws.Range("I122").Value = 424.5
$ws.Range("J122").Value = 272.5
$ws.Range("K122").Value = 3820.5
$ws.Range("L122").Value = 2452.5
$ws.Range("M122").Value = -1370.5
$ws.Range("N122").Value = -7352.5

$ws.Range("H132").Value = 2562.889
$ws.Range("I132").Value = 2513.8333
$ws.Range("J132").Value = 2661
$ws.Range("K132").Value = 22624.4997
$ws.Range("L132").Value = 23949
$ws.Range("M132").Value = -20094.4997
$ws.Range("N132").Value = -29009

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 15000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 15000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 15000
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -15710

$ws.Range("H80").Value = 5548.8945
$ws.Range("I80").Value = 3621
$ws.Range("K80").Value = 3621
$ws.Range("M80").Value = -2623

$ws.Range("H83").Value = 5548.8945
$ws.Range("I83").Value = 3621
$ws.Range("K83").Value = 18105
$ws.Range("M83").Value = -13113

$ws.Range("H97").Value = 856.38464
$ws.Range("I97").Value = 495.83334
$ws.Range("J97").Value = 1165.4286
$ws.Range("K97").Value = 495.83334
$ws.Range("L97").Value = 1165.4286
$ws.Range("M97").Value = 0.1666599999999789
$ws.Range("N97").Value = -2157.4286

$ws.Range("H126").Value = 5109.6
$ws.Range("I126").Value = 5102.6665
$ws.Range("J126").Value = 5120
$ws.Range("K126").Value = 15307.9995
$ws.Range("L126").Value = 15360
$ws.Range("M126").Value = -12837.9995
$ws.Range("N126").Value = -20300

$ws.Range("H132").Value = 3979.6
$ws.Range("I132").Value = 1974.75
$ws.Range("K132").Value = 5924.25
$ws.Range("M132").Value = -3394.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5458.963
$ws.Range("I40").Value = 4768.136
$ws.Range("K40").Value = 4768.136
$ws.Range("M40").Value = -4632.136

$ws.Range("H46").Value = 2510.5557
$ws.Range("I46").Value = 922.3333
$ws.Range("J46").Value = 3781.1333
$ws.Range("K46").Value = 922.3333
$ws.Range("L46").Value = 3781.1333
$ws.Range("M46").Value = -734.3333
$ws.Range("N46").Value = -4157.1333

$ws.Range("H68").Value = 4975
$ws.Range("I68").Value = 4975
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 4975
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -4226
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 4975
$ws.Range("I71").Value = 4975
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 24875
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -21131
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1349.3889
$ws.Range("J81").Value = 1056.7142
$ws.Range("L81").Value = 2113.4284
$ws.Range("N81").Value = -4235.4284

$ws.Range("H84").Value = 1349.3889
$ws.Range("J84").Value = 1056.7142
$ws.Range("L84").Value = 10567.142
$ws.Range("N84").Value = -21175.142

$ws.Range("H122").Value = 1791.7959
$ws.Range("I122").Value = 1815.1163
$ws.Range("J122").Value = 1624.6666
$ws.Range("K122").Value = 5445.3489
$ws.Range("L122").Value = 4873.9998
$ws.Range("M122").Value = -2995.3489
$ws.Range("N122").Value = -9773.9998

$ws.Range("H132").Value = 2945.2927
$ws.Range("I132").Value = 1955.2941
$ws.Range("J132").Value = 7753.857
$ws.Range("K132").Value = 5865.8823
$ws.Range("L132").Value = 23261.571
$ws.Range("M132").Value = -3335.8823
$ws.Range("N132").Value = -28321.571
